$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.643452
$ws.Range("H2").Value = 1.930356
$ws.Range("I2").Value = 0.2431136893481813
$ws.Range("J2").Value = 0.2431136893481813
$ws.Range("M2").Value = 1.309671333333333
$ws.Range("N2").Value = 3.929014
$ws.Range("O2").Value = 0.05806924226264097
$ws.Range("P2").Value = 0.05806924226264098
$ws.Range("Q2").Value = 0.8427106387760001
$ws.Range("R2").Value = 7.584395748984
$ws.Range("S2").Value = 0.01411742772412398
$ws.Range("T2").Value = 0.01411742772412398
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.643452
$ws.Range("H3").Value = 1.930356
$ws.Range("I3").Value = 0.2431136893481813
$ws.Range("J3").Value = 0.2431136893481813
$ws.Range("O3").Value = 0.3245116581089107
$ws.Range("P3").Value = 0.3245116581089107
$ws.Range("Q3").Value = 4.709367920772
$ws.Range("R3").Value = 42.384311286948
$ws.Range("S3").Value = 0.07889322643935293
$ws.Range("T3").Value = 0.07889322643935293
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.643452
$ws.Range("H4").Value = 1.930356
$ws.Range("I4").Value = 0.2431136893481813
$ws.Range("J4").Value = 0.2431136893481813
$ws.Range("M4").Value = 4.657910333333334
$ws.Range("N4").Value = 13.973731
$ws.Range("O4").Value = 0.2065261082683789
$ws.Range("P4").Value = 0.2065261082683789
$ws.Range("Q4").Value = 2.997141719804
$ws.Range("R4").Value = 26.974275478236
$ws.Range("S4").Value = 0.05020932412784751
$ws.Range("T4").Value = 0.05020932412784751
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.643452
$ws.Range("H5").Value = 1.930356
$ws.Range("I5").Value = 0.2431136893481813
$ws.Range("J5").Value = 0.2431136893481813
$ws.Range("M5").Value = 9.267122333333333
$ws.Range("N5").Value = 27.801367
$ws.Range("O5").Value = 0.4108929913600695
$ws.Range("P5").Value = 0.4108929913600695
$ws.Range("Q5").Value = 5.962948399628
$ws.Range("R5").Value = 53.666535596652
$ws.Range("S5").Value = 0.09989371105685686
$ws.Range("T5").Value = 0.09989371105685686
$ws.Range("I6").Value = 0.2185191514957488
$ws.Range("J6").Value = 0.2185191514957488
$ws.Range("M6").Value = 1.309671333333333
$ws.Range("N6").Value = 3.929014
$ws.Range("O6").Value = 0.05806924226264097
$ws.Range("P6").Value = 0.05806924226264098
$ws.Range("Q6").Value = 0.7574580198897779
$ws.Range("R6").Value = 6.817122179008001
$ws.Range("S6").Value = 0.01268924154723338
$ws.Range("T6").Value = 0.01268924154723338
$ws.Range("I7").Value = 0.2185191514957488
$ws.Range("J7").Value = 0.2185191514957488
$ws.Range("O7").Value = 0.3245116581089107
$ws.Range("P7").Value = 0.3245116581089107
$ws.Range("S7").Value = 0.07091201218043769
$ws.Range("T7").Value = 0.0709120121804377
$ws.Range("I8").Value = 0.2185191514957488
$ws.Range("J8").Value = 0.2185191514957488
$ws.Range("M8").Value = 4.657910333333334
$ws.Range("N8").Value = 13.973731
$ws.Range("O8").Value = 0.2065261082683789
$ws.Range("P8").Value = 0.2065261082683789
$ws.Range("Q8").Value = 2.693936599292445
$ws.Range("R8").Value = 24.245429393632
$ws.Range("S8").Value = 0.04512990994052529
$ws.Range("T8").Value = 0.0451299099405253
$ws.Range("I9").Value = 0.2185191514957488
$ws.Range("J9").Value = 0.2185191514957488
$ws.Range("M9").Value = 9.267122333333333
$ws.Range("N9").Value = 27.801367
$ws.Range("O9").Value = 0.4108929913600695
$ws.Range("P9").Value = 0.4108929913600695
$ws.Range("Q9").Value = 5.359708160380444
$ws.Range("R9").Value = 48.237373443424
$ws.Range("S9").Value = 0.08978798782755242
$ws.Range("T9").Value = 0.08978798782755242
$ws.Range("G10").Value = 0.2588786666666666
$ws.Range("H10").Value = 0.776636
$ws.Range("I10").Value = 0.09781141055878506
$ws.Range("J10").Value = 0.09781141055878506
$ws.Range("M10").Value = 1.309671333333333
$ws.Range("N10").Value = 3.929014
$ws.Range("O10").Value = 0.05806924226264097
$ws.Range("P10").Value = 0.05806924226264098
$ws.Range("Q10").Value = 0.3390459685448889
$ws.Range("R10").Value = 3.051413716904
$ws.Range("S10").Value = 0.005679834495788729
$ws.Range("T10").Value = 0.005679834495788729
$ws.Range("G11").Value = 0.2588786666666666
$ws.Range("H11").Value = 0.776636
$ws.Range("I11").Value = 0.09781141055878506
$ws.Range("J11").Value = 0.09781141055878506
$ws.Range("O11").Value = 0.3245116581089107
$ws.Range("P11").Value = 0.3245116581089107
$ws.Range("Q11").Value = 1.894709921132
$ws.Range("R11").Value = 17.052389290188
$ws.Range("S11").Value = 0.03174094302240275
$ws.Range("T11").Value = 0.03174094302240276
$ws.Range("G12").Value = 0.2588786666666666
$ws.Range("H12").Value = 0.776636
$ws.Range("I12").Value = 0.09781141055878506
$ws.Range("J12").Value = 0.09781141055878506
$ws.Range("M12").Value = 4.657910333333334
$ws.Range("N12").Value = 13.973731
$ws.Range("O12").Value = 0.2065261082683789
$ws.Range("P12").Value = 0.2065261082683789
$ws.Range("Q12").Value = 1.205833616546222
$ws.Range("R12").Value = 10.852502548916
$ws.Range("S12").Value = 0.0202006099669465
$ws.Range("T12").Value = 0.0202006099669465
$ws.Range("G13").Value = 0.2588786666666666
$ws.Range("H13").Value = 0.776636
$ws.Range("I13").Value = 0.09781141055878506
$ws.Range("J13").Value = 0.09781141055878506
$ws.Range("M13").Value = 9.267122333333333
$ws.Range("N13").Value = 27.801367
$ws.Range("O13").Value = 0.4108929913600695
$ws.Range("P13").Value = 0.4108929913600695
$ws.Range("Q13").Value = 2.399060273490222
$ws.Range("R13").Value = 21.591542461412
$ws.Range("S13").Value = 0.04019002307364708
$ws.Range("T13").Value = 0.04019002307364708
$ws.Range("G14").Value = 1.166024333333333
$ws.Range("H14").Value = 3.498073
$ws.Range("I14").Value = 0.4405557485972849
$ws.Range("J14").Value = 0.4405557485972849
$ws.Range("M14").Value = 1.309671333333333
$ws.Range("N14").Value = 3.929014
$ws.Range("O14").Value = 0.05806924226264097
$ws.Range("P14").Value = 0.05806924226264098
$ws.Range("Q14").Value = 1.527108643335778
$ws.Range("R14").Value = 13.743977790022
$ws.Range("S14").Value = 0.02558273849549489
$ws.Range("T14").Value = 0.02558273849549489
$ws.Range("G15").Value = 1.166024333333333
$ws.Range("H15").Value = 3.498073
$ws.Range("I15").Value = 0.4405557485972849
$ws.Range("J15").Value = 0.4405557485972849
$ws.Range("O15").Value = 0.3245116581089107
$ws.Range("P15").Value = 0.3245116581089107
$ws.Range("Q15").Value = 8.534028319500999
$ws.Range("R15").Value = 76.806254875509
$ws.Range("S15").Value = 0.1429654764667173
$ws.Range("T15").Value = 0.1429654764667173
$ws.Range("G16").Value = 1.166024333333333
$ws.Range("H16").Value = 3.498073
$ws.Range("I16").Value = 0.4405557485972849
$ws.Range("J16").Value = 0.4405557485972849
$ws.Range("M16").Value = 4.657910333333334
$ws.Range("N16").Value = 13.973731
$ws.Range("O16").Value = 0.2065261082683789
$ws.Range("P16").Value = 0.2065261082683789
$ws.Range("Q16").Value = 5.431236791151445
$ws.Range("R16").Value = 48.881131120363
$ws.Range("S16").Value = 0.09098626423305957
$ws.Range("T16").Value = 0.09098626423305957
$ws.Range("G17").Value = 1.166024333333333
$ws.Range("H17").Value = 3.498073
$ws.Range("I17").Value = 0.4405557485972849
$ws.Range("J17").Value = 0.4405557485972849
$ws.Range("M17").Value = 9.267122333333333
$ws.Range("N17").Value = 27.801367
$ws.Range("O17").Value = 0.4108929913600695
$ws.Range("P17").Value = 0.4108929913600695
$ws.Range("Q17").Value = 10.80569014064344
$ws.Range("R17").Value = 97.25121126579099
$ws.Range("S17").Value = 0.1810212694020131
$ws.Range("T17").Value = 0.1810212694020131
